# The commit inserts one new daily price-report row into the "Arándano (blue)"
# sheet at position 39 (pushing every existing row below it down by one), and
# fills the new row with a fresh sample of the same shape as the surrounding
# rows. The sheet's used range therefore grows from A1:T130 to A1:T131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 39; Excel shifts rows 39..130 down to 40..131 and
# extends the sheet dimension automatically (carries the date-format style
# from the surrounding rows onto the new D39 cell as well).
$ws.Rows("39:39").Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Macroferia Regional de Talca"
$ws.Range("C39").Value = "Maule"
$ws.Range("D39").Value = 45002
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100101
$ws.Range("H39").Value = "Berries"
$ws.Range("I39").Value = 100101001
$ws.Range("J39").Value = "Arándano (blue)"
$ws.Range("K39").Value = "Sin especificar"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 30
$ws.Range("N39").Value = 4000
$ws.Range("O39").Value = 4000
$ws.Range("P39").Value = 4000
$ws.Range("Q39").Value = "`$/bandeja 2 kilos"
$ws.Range("R39").Value = "Provincia de Curicó"
$ws.Range("S39").Value = 2000
$ws.Range("T39").Value = 2
